# 25th commit Final 001-017 Updated
#
# This edit:
#  1. Updates the "V8" user on Basic_Information_Release (D2) to reference
#     the newly-added "v7_automation_v8_user [v7_automation_v8_user]" string.
#  2. Updates tc_003's V8 / Project Write Access columns (D2, E2) from
#     "Sahoo Sukriti [suksaho]" to "Avinash Ankush [ANAVINA]".
#  3. Moves the active selection around a couple of sheets and re-activates
#     tc_003 as the selected/active tab (was tc_004).
#  4. Narrows the "V8" column on Basic_Information_Release.

$wb = $excel.ActiveWorkbook

# --- Basic_Information_Release -------------------------------------------------
$wsBasic = $wb.Worksheets.Item("Basic_Information_Release")
$wsBasic.Range("D2").Value = "v7_automation_v8_user [v7_automation_v8_user]"
$wsBasic.Columns("D").ColumnWidth = 20
$wsBasic.Activate() | Out-Null
$wsBasic.Range("D11").Select() | Out-Null

# --- tc_004 ----------------------------------------------------------------
$wsTc004 = $wb.Worksheets.Item("tc_004")
$wsTc004.Activate() | Out-Null
$wsTc004.Range("C6").Select() | Out-Null

# --- tc_003 (becomes the active / tab-selected sheet) -----------------------
$wsTc003 = $wb.Worksheets.Item("tc_003")
$wsTc003.Range("D2").Value = "Avinash Ankush [ANAVINA]"
$wsTc003.Range("E2").Value = "Avinash Ankush [ANAVINA]"
$wsTc003.Activate() | Out-Null
$wsTc003.Range("E2").Select() | Out-Null
